# Apply updated distributive results figures (Arthen's reduction-in-dividends
# implementation + small adjustments) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Regime Atual"
$ws.Range("B2").Value = 0.6185935638155903
$ws.Range("C2").Value = 0.1401127577307753
$ws.Range("D2").Value = 0.5492320602264944
$ws.Range("E2").Value = 0.2478619132032327
$ws.Range("F2").Value = 334.6955778851641

# Row 3 - "Nova Proposta"
$ws.Range("B3").Value = 0.6165899422660595
$ws.Range("C3").Value = 0.1406696751015247
$ws.Range("D3").Value = 0.5466961269590717
$ws.Range("E3").Value = 0.2421434918316321
$ws.Range("F3").Value = 360.6340583975178
$ws.Range("G3").Value = 25.93848051235375

# Row 4 - "Nova c/ Aliq. Máxima"
$ws.Range("B4").Value = 0.6152994338053641
$ws.Range("C4").Value = 0.1411481392112925
$ws.Range("D4").Value = 0.5451542904977854
$ws.Range("E4").Value = 0.2395228084753867
$ws.Range("F4").Value = 382.5895231548401
$ws.Range("G4").Value = 47.89394526967601
